# Fixed harvester column in rnaSamples -- holly added S.GISH to harvester in bioSamples
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column B (harvester) like the author did while reviewing/fixing it
$ws.Columns("B").Select()

# Replace the harvester value for every data row with the corrected entry
$ws.Range("B2:B25").Value = "S.GISH"
